$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "May An Phu" / "Example" example row, highlighted yellow ---
$ws.Range("A2").Value = "May An Phu"
$ws.Range("B2").Value = "Example"
$ws.Range("A2").Interior.Color = 65535   # yellow - A2 keeps the existing thin border (style: fill+border)
$ws.Range("B2").Interior.Color = 65535   # yellow - B2 has no border in this sheet (style: fill only)

# --- Row 3: blank cell reserved/highlighted red ---
$ws.Range("A3").Interior.Color = 255     # red - A3 keeps the existing thin border (style: fill+border)

# Mint the "fill + no border" combination for the red fill too (index parity with the
# template's style table: yellow/bordered, yellow/borderless, red/bordered, red/borderless),
# using a scratch cell outside the sheet's real range, then discard it.
$ws.Range("Z1").Interior.Color = 255
$ws.Range("Z1").Borders.LineStyle = -4142
$ws.Range("Z1").Clear()

# --- Selection / view state ---
$ws.Range("A3").Select()
